$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2672.2  # ALC!H17: 2587.4546 -> 2672.2
$ws.Cells.Item(17, 10).Value = 2858  # ALC!J17: 2746.2 -> 2858
$ws.Cells.Item(17, 12).Value = 8574  # ALC!L17: 8238.599999999999 -> 8574
$ws.Cells.Item(17, 14).Value = -8910  # ALC!N17: -8574.599999999999 -> -8910

$ws.Cells.Item(80, 8).Value = 920  # ALC!H80: 834.0714 -> 920
$ws.Cells.Item(80, 9).Value = 1225  # ALC!I80: 944 -> 1225
$ws.Cells.Item(80, 10).Value = 852.2222  # ALC!J80: 804.0909 -> 852.2222
$ws.Cells.Item(80, 11).Value = 3675  # ALC!K80: 2832 -> 3675
$ws.Cells.Item(80, 12).Value = 2556.6666  # ALC!L80: 2412.2727 -> 2556.6666
$ws.Cells.Item(80, 13).Value = -2677  # ALC!M80: -1834 -> -2677
$ws.Cells.Item(80, 14).Value = -4552.6666  # ALC!N80: -4408.2727 -> -4552.6666

$ws.Cells.Item(83, 8).Value = 920  # ALC!H83: 834.0714 -> 920
$ws.Cells.Item(83, 9).Value = 1225  # ALC!I83: 944 -> 1225
$ws.Cells.Item(83, 10).Value = 852.2222  # ALC!J83: 804.0909 -> 852.2222
$ws.Cells.Item(83, 11).Value = 11025  # ALC!K83: 8496 -> 11025
$ws.Cells.Item(83, 12).Value = 7669.999800000001  # ALC!L83: 7236.8181 -> 7669.999800000001
$ws.Cells.Item(83, 13).Value = -6033  # ALC!M83: -3504 -> -6033
$ws.Cells.Item(83, 14).Value = -17653.9998  # ALC!N83: -17220.8181 -> -17653.9998

$ws.Cells.Item(111, 8).Value = 1614.6  # ALC!H111: 1325.3334 -> 1614.6
$ws.Cells.Item(111, 9).Value = 1484.3334  # ALC!I111: 1238 -> 1484.3334
$ws.Cells.Item(111, 10).Value = 1810  # ALC!J111: 1500 -> 1810
$ws.Cells.Item(111, 11).Value = 4453.0002  # ALC!K111: 3714 -> 4453.0002
$ws.Cells.Item(111, 12).Value = 5430  # ALC!L111: 4500 -> 5430
$ws.Cells.Item(111, 13).Value = -1386.0002  # ALC!M111: -647 -> -1386.0002
$ws.Cells.Item(111, 14).Value = -11564  # ALC!N111: -10634 -> -11564

$ws.Cells.Item(127, 8).Value = 499  # ALC!H127: 0 -> 499
$ws.Cells.Item(127, 9).Value = 499  # ALC!I127: 0 -> 499
$ws.Cells.Item(127, 11).Value = 1497  # ALC!K127: 0 -> 1497
$ws.Cells.Item(127, 13).Value = 3463  # ALC!M127: (new) -> 3463

$ws.Cells.Item(132, 8).Value = 10611.75  # ALC!H132: 10620.35 -> 10611.75
$ws.Cells.Item(132, 9).Value = 10611.75  # ALC!I132: 10620.35 -> 10611.75
$ws.Cells.Item(132, 11).Value = 31835.25  # ALC!K132: 31861.05 -> 31835.25
$ws.Cells.Item(132, 13).Value = -29305.25  # ALC!M132: -29331.05 -> -29305.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 6628.952  # ARM!H2: 6367.9546 -> 6628.952
$ws.Cells.Item(2, 9).Value = 3089.9473  # ARM!I2: 2979.8 -> 3089.9473
$ws.Cells.Item(2, 11).Value = 3089.9473  # ARM!K2: 2979.8 -> 3089.9473
$ws.Cells.Item(2, 13).Value = -2976.9473  # ARM!M2: -2866.8 -> -2976.9473

$ws.Cells.Item(38, 8).Value = 3754.5  # ARM!H38: 4673 -> 3754.5
$ws.Cells.Item(38, 9).Value = 3754.5  # ARM!I38: 4673 -> 3754.5
$ws.Cells.Item(38, 11).Value = 3754.5  # ARM!K38: 4673 -> 3754.5
$ws.Cells.Item(38, 13).Value = -3287.5  # ARM!M38: -4206 -> -3287.5

$ws.Cells.Item(45, 8).Value = 2747.8823  # ARM!H45: 2734.5264 -> 2747.8823
$ws.Cells.Item(45, 9).Value = 1842.2727  # ARM!I45: 1854.4166 -> 1842.2727
$ws.Cells.Item(45, 10).Value = 4408.1665  # ARM!J45: 4243.2856 -> 4408.1665
$ws.Cells.Item(45, 11).Value = 1842.2727  # ARM!K45: 1854.4166 -> 1842.2727
$ws.Cells.Item(45, 12).Value = 4408.1665  # ARM!L45: 4243.2856 -> 4408.1665
$ws.Cells.Item(45, 13).Value = -1465.2727  # ARM!M45: -1477.4166 -> -1465.2727
$ws.Cells.Item(45, 14).Value = -5162.1665  # ARM!N45: -4997.2856 -> -5162.1665

$ws.Cells.Item(74, 8).Value = 3050.5  # ARM!H74: 3600.9 -> 3050.5
$ws.Cells.Item(74, 9).Value = 2509.6365  # ARM!I74: 3001 -> 2509.6365
$ws.Cells.Item(74, 11).Value = 2509.6365  # ARM!K74: 3001 -> 2509.6365
$ws.Cells.Item(74, 13).Value = -1635.6365  # ARM!M74: -2127 -> -1635.6365

$ws.Cells.Item(77, 8).Value = 3050.5  # ARM!H77: 3600.9 -> 3050.5
$ws.Cells.Item(77, 9).Value = 2509.6365  # ARM!I77: 3001 -> 2509.6365
$ws.Cells.Item(77, 11).Value = 12548.1825  # ARM!K77: 15005 -> 12548.1825
$ws.Cells.Item(77, 13).Value = -8180.182500000001  # ARM!M77: -10637 -> -8180.182500000001

$ws.Cells.Item(97, 8).Value = 3477.75  # ARM!H97: 3141.2 -> 3477.75
$ws.Cells.Item(97, 9).Value = 2563  # ARM!I97: 2578 -> 2563
$ws.Cells.Item(97, 10).Value = 6222  # ARM!J97: 3986 -> 6222
$ws.Cells.Item(97, 11).Value = 2563  # ARM!K97: 2578 -> 2563
$ws.Cells.Item(97, 12).Value = 6222  # ARM!L97: 3986 -> 6222
$ws.Cells.Item(97, 13).Value = -2067  # ARM!M97: -2082 -> -2067
$ws.Cells.Item(97, 14).Value = -7214  # ARM!N97: -4978 -> -7214

$ws.Cells.Item(116, 8).Value = 6628.952  # ARM!H116: 6367.9546 -> 6628.952
$ws.Cells.Item(116, 9).Value = 3089.9473  # ARM!I116: 2979.8 -> 3089.9473
$ws.Cells.Item(116, 11).Value = 3089.9473  # ARM!K116: 2979.8 -> 3089.9473
$ws.Cells.Item(116, 13).Value = -795.9472999999998  # ARM!M116: -685.8000000000002 -> -795.9472999999998

$ws.Cells.Item(132, 8).Value = 2249.1667  # ARM!H132: 2532.6667 -> 2249.1667
$ws.Cells.Item(132, 9).Value = 2249.1667  # ARM!I132: 2532.6667 -> 2249.1667
$ws.Cells.Item(132, 11).Value = 6747.500100000001  # ARM!K132: 7598.000100000001 -> 6747.500100000001
$ws.Cells.Item(132, 13).Value = -4217.500100000001  # ARM!M132: -5068.000100000001 -> -4217.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 6628.952  # BSM!H3: 6367.9546 -> 6628.952
$ws.Cells.Item(3, 9).Value = 3089.9473  # BSM!I3: 2979.8 -> 3089.9473
$ws.Cells.Item(3, 11).Value = 3089.9473  # BSM!K3: 2979.8 -> 3089.9473
$ws.Cells.Item(3, 13).Value = -2975.9473  # BSM!M3: -2865.8 -> -2975.9473

$ws.Cells.Item(94, 8).Value = 825  # BSM!H94: 875 -> 825
$ws.Cells.Item(94, 9).Value = 800  # BSM!I94: 860 -> 800
$ws.Cells.Item(94, 10).Value = 862.5  # BSM!J94: 900 -> 862.5
$ws.Cells.Item(94, 11).Value = 800  # BSM!K94: 860 -> 800
$ws.Cells.Item(94, 12).Value = 862.5  # BSM!L94: 900 -> 862.5
$ws.Cells.Item(94, 13).Value = -349  # BSM!M94: -409 -> -349
$ws.Cells.Item(94, 14).Value = -1764.5  # BSM!N94: -1802 -> -1764.5

$ws.Cells.Item(99, 8).Value = 4855.8887  # BSM!H99: 5969.6665 -> 4855.8887
$ws.Cells.Item(99, 9).Value = 4868  # BSM!I99: 6454.5 -> 4868
$ws.Cells.Item(99, 10).Value = 4831.6665  # BSM!J99: 5000 -> 4831.6665
$ws.Cells.Item(99, 11).Value = 4868  # BSM!K99: 6454.5 -> 4868
$ws.Cells.Item(99, 12).Value = 4831.6665  # BSM!L99: 5000 -> 4831.6665
$ws.Cells.Item(99, 13).Value = -3370  # BSM!M99: -4956.5 -> -3370
$ws.Cells.Item(99, 14).Value = -7827.6665  # BSM!N99: -7996 -> -7827.6665

$ws.Cells.Item(134, 8).Value = 1899.1904  # BSM!H134: 1937.0952 -> 1899.1904
$ws.Cells.Item(134, 9).Value = 1371.125  # BSM!I134: 1382.5333 -> 1371.125
$ws.Cells.Item(134, 10).Value = 3589  # BSM!J134: 3323.5 -> 3589
$ws.Cells.Item(134, 11).Value = 4113.375  # BSM!K134: 4147.5999 -> 4113.375
$ws.Cells.Item(134, 12).Value = 10767  # BSM!L134: 9970.5 -> 10767
$ws.Cells.Item(134, 13).Value = -1578.375  # BSM!M134: -1612.5999 -> -1578.375
$ws.Cells.Item(134, 14).Value = -15837  # BSM!N134: -15040.5 -> -15837

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(32, 8).Value = 1535.2858  # CRP!H32: 1417.8334 -> 1535.2858
$ws.Cells.Item(32, 9).Value = 1416.1666  # CRP!I32: 1417.8334 -> 1416.1666
$ws.Cells.Item(32, 10).Value = 2250  # CRP!J32: 0 -> 2250
$ws.Cells.Item(32, 11).Value = 1416.1666  # CRP!K32: 1417.8334 -> 1416.1666
$ws.Cells.Item(32, 12).Value = 2250  # CRP!L32: 0 -> 2250
$ws.Cells.Item(32, 13).Value = -1100.1666  # CRP!M32: -1101.8334 -> -1100.1666
$ws.Cells.Item(32, 14).Value = -2882  # CRP!N32: (new) -> -2882

$ws.Cells.Item(53, 8).Value = 54842  # CRP!H53: 55000 -> 54842
$ws.Cells.Item(53, 10).Value = 54842  # CRP!J53: 55000 -> 54842
$ws.Cells.Item(53, 12).Value = 54842  # CRP!L53: 55000 -> 54842
$ws.Cells.Item(53, 14).Value = -56056  # CRP!N53: -56214 -> -56056

$ws.Cells.Item(111, 8).Value = 62836  # CRP!H111: 0 -> 62836
$ws.Cells.Item(111, 10).Value = 62836  # CRP!J111: 0 -> 62836
$ws.Cells.Item(111, 12).Value = 62836  # CRP!L111: 0 -> 62836
$ws.Cells.Item(111, 14).Value = -71016  # CRP!N111: (new) -> -71016

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 2746.25  # CUL!H141: 2997.5 -> 2746.25
$ws.Cells.Item(141, 9).Value = 2746.25  # CUL!I141: 2997.5 -> 2746.25
$ws.Cells.Item(141, 11).Value = 8238.75  # CUL!K141: 8992.5 -> 8238.75
$ws.Cells.Item(141, 13).Value = -3058.75  # CUL!M141: -3812.5 -> -3058.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 24000  # GSM!H26: 20508.4 -> 24000
$ws.Cells.Item(26, 10).Value = 23500  # GSM!J26: 19385.5 -> 23500
$ws.Cells.Item(26, 12).Value = 23500  # GSM!L26: 19385.5 -> 23500
$ws.Cells.Item(26, 14).Value = -24060  # GSM!N26: -19945.5 -> -24060

$ws.Cells.Item(50, 8).Value = 24000  # GSM!H50: 20508.4 -> 24000
$ws.Cells.Item(50, 10).Value = 23500  # GSM!J50: 19385.5 -> 23500
$ws.Cells.Item(50, 12).Value = 23500  # GSM!L50: 19385.5 -> 23500
$ws.Cells.Item(50, 14).Value = -24496  # GSM!N50: -20381.5 -> -24496

$ws.Cells.Item(97, 8).Value = 0  # GSM!H97: 1449 -> 0
$ws.Cells.Item(97, 9).Value = 0  # GSM!I97: 1449 -> 0
$ws.Cells.Item(97, 11).Value = 0  # GSM!K97: 1449 -> 0
$ws.Cells.Item(97, 13).Value = $null  # GSM!M97: -> (removed)

$ws.Cells.Item(132, 8).Value = 4868.241  # GSM!H132: 4898.8125 -> 4868.241
$ws.Cells.Item(132, 9).Value = 4320.35  # GSM!I132: 4412.625 -> 4320.35
$ws.Cells.Item(132, 10).Value = 6085.778  # GSM!J132: 6357.375 -> 6085.778
$ws.Cells.Item(132, 11).Value = 12961.05  # GSM!K132: 13237.875 -> 12961.05
$ws.Cells.Item(132, 12).Value = 18257.334  # GSM!L132: 19072.125 -> 18257.334
$ws.Cells.Item(132, 13).Value = -10431.05  # GSM!M132: -10707.875 -> -10431.05
$ws.Cells.Item(132, 14).Value = -23317.334  # GSM!N132: -24132.125 -> -23317.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5792.154  # LTW!H7: 6108.25 -> 5792.154
$ws.Cells.Item(7, 9).Value = 4181.875  # LTW!I7: 4493.7144 -> 4181.875
$ws.Cells.Item(7, 11).Value = 4181.875  # LTW!K7: 4493.7144 -> 4181.875
$ws.Cells.Item(7, 13).Value = -4069.875  # LTW!M7: -4381.7144 -> -4069.875

$ws.Cells.Item(61, 8).Value = 4517.091  # LTW!H61: 4818.8 -> 4517.091
$ws.Cells.Item(61, 9).Value = 1615  # LTW!I61: 1638 -> 1615
$ws.Cells.Item(61, 11).Value = 1615  # LTW!K61: 1638 -> 1615
$ws.Cells.Item(61, 13).Value = -1413  # LTW!M61: -1436 -> -1413

$ws.Cells.Item(113, 8).Value = 4517.091  # LTW!H113: 4818.8 -> 4517.091
$ws.Cells.Item(113, 9).Value = 1615  # LTW!I113: 1638 -> 1615
$ws.Cells.Item(113, 11).Value = 1615  # LTW!K113: 1638 -> 1615
$ws.Cells.Item(113, 13).Value = 555  # LTW!M113: 532 -> 555

$ws.Cells.Item(126, 8).Value = 5792.154  # LTW!H126: 6108.25 -> 5792.154
$ws.Cells.Item(126, 9).Value = 4181.875  # LTW!I126: 4493.7144 -> 4181.875
$ws.Cells.Item(126, 11).Value = 12545.625  # LTW!K126: 13481.1432 -> 12545.625
$ws.Cells.Item(126, 13).Value = -10075.625  # LTW!M126: -11011.1432 -> -10075.625

$ws.Cells.Item(132, 8).Value = 2500  # LTW!H132: 3062.3 -> 2500
$ws.Cells.Item(132, 9).Value = 2500  # LTW!I132: 3015.5 -> 2500
$ws.Cells.Item(132, 10).Value = 0  # LTW!J132: 3249.5 -> 0
$ws.Cells.Item(132, 11).Value = 7500  # LTW!K132: 9046.5 -> 7500
$ws.Cells.Item(132, 12).Value = 0  # LTW!L132: 9748.5 -> 0
$ws.Cells.Item(132, 13).Value = $null  # LTW!M132: -6516.5 -> (removed)
$ws.Cells.Item(132, 14).Value = $null  # LTW!N132: -> (removed)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2212.4167  # WVR!H132: 2657.0952 -> 2212.4167
$ws.Cells.Item(132, 9).Value = 1033  # WVR!I132: 1429.3077 -> 1033
$ws.Cells.Item(132, 10).Value = 4178.1113  # WVR!J132: 4652.25 -> 4178.1113
$ws.Cells.Item(132, 11).Value = 3099  # WVR!K132: 4287.9231 -> 3099
$ws.Cells.Item(132, 12).Value = 12534.3339  # WVR!L132: 13956.75 -> 12534.3339
$ws.Cells.Item(132, 13).Value = -569  # WVR!M132: -1757.9231 -> -569
$ws.Cells.Item(132, 14).Value = -17594.3339  # WVR!N132: -19016.75 -> -17594.3339

$ws.Cells.Item(138, 8).Value = 0  # WVR!H138: 100000 -> 0
$ws.Cells.Item(138, 10).Value = 0  # WVR!J138: 100000 -> 0
$ws.Cells.Item(138, 12).Value = 0  # WVR!L138: 100000 -> 0
$ws.Cells.Item(138, 14).Value = 0  # WVR!N138: -110280 -> 0
